$d = $word.ActiveDocument

# The page footer (blank line, "Ver no Jupiter ..." and the "(c) 2020 ..."
# copyright notice) that used to trail the "LOQ4237: ... (Requisito fraco)"
# requirement line was dropped from this rebuilt page. Locate it by its
# text and delete the whole block, leaving the following blank paragraph
# (and the page-break paragraph after it) untouched.

$loq = $d.Content
$null = $loq.Find.Execute("LOQ4237: Projeto Integrado de Engenharia de Produção II (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deleteStart = $loq.Paragraphs(1).Range.End

$copyright = $d.Content
$null = $copyright.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deleteEnd = $copyright.Paragraphs(1).Range.End

$d.Range($deleteStart, $deleteEnd).Delete()
